$p = $ppt.ActivePresentation

# Slide 13 ("CAST & CONVERT" / ISNULL overview slide): the closing phrase of
# the first bullet changes from "преобразуване между дати" (conversion
# between dates) to "преобразуване между данни" (conversion between data),
# which PowerPoint records as the edited word-group being split into its
# own run.
$slide = $p.Slides.Item(13)
$shape = $slide.Shapes.Item(8)
$tr = $shape.TextFrame.TextRange
[void]$tr.Replace("между дати", "между данни")
